$wb = $excel.ActiveWorkbook

# The original (and only) sheet holds all the diary rows -> becomes "part1"
$wb.Worksheets.Item(1).Name = "part1"

# Add two new, empty sheets for the other parts. Worksheets.Add() always
# inserts right before the currently active sheet and makes the new sheet
# active, so adding "part2" then "part3" in that order yields tab order
# part3, part2, part1 - matching the target workbook.
$wb.Worksheets.Add().Name = "part2"
$wb.Worksheets.Add().Name = "part3"

# Give the new sheets the same header row + hour-total formula as part1.
# Re-look-up each sheet by name every time instead of reusing a captured
# object reference - references captured before other sheets were
# inserted/removed end up stale and silently touch the wrong sheet.
foreach ($name in @("part2", "part3")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1").Value = "pvm"
    $ws.Range("B1").Value = "time (min)"
    $ws.Range("C1").Value = "comment"
    $ws.Range("F1").Value = "time (hour)"
    $ws.Range("F2").Formula = "=SUM(B2:B38)/60"
}

# part1 no longer needs the "part" column (always 1) - drop column C and
# let "comment" (old D) and "time (hour)" (old G) shift left.
$wb.Worksheets.Item("part1").Columns.Item(3).Delete()

# Record today's diary submission as a new row.
$wb.Worksheets.Item("part1").Range("A11").Value = 211101
$wb.Worksheets.Item("part1").Range("B11").Value = 10
$wb.Worksheets.Item("part1").Range("C11").Value = "Submit and edit diary"

# part1 is the tab that should remain selected/active.
$wb.Worksheets.Item("part1").Select()
